$d = $word.ActiveDocument
$apos = [char]0x2019

function Find-ParaIndex($doc, $wanted) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx++
        $t = $p.Range.Text.TrimEnd([char]0x0D, [char]0x07)
        if ($t -eq $wanted) {
            return $idx
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Merge the split run in the "Gps technology..." bullet into one run:
#    " technology can be implemented" + " as a potential use"
#       -> " technology can be implemented as a potential use"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " technology can be implemented as a potential use",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " technology can be implemented as a potential use", 2)

# ------------------------------------------------------------------
# 2) Move "Parking" from the numId=1 list (right after "Dining Services")
#    down into the numId=3 list, landing between "Interactive map of
#    campus itself as a whole" and "3-Dimensional campus overview".
#    We delete the old bullet and grow the numId=3 list instead of
#    reusing the cut paragraph, so the new bullet naturally inherits
#    numId=3 (matching the target numbering) rather than dragging its
#    old numId=1 formatting along.
# ------------------------------------------------------------------
$parkIdx = Find-ParaIndex $d "Parking"
$parkingOld = $d.Paragraphs($parkIdx)
$parkingOld.Range.Delete()

$afterIdx = Find-ParaIndex $d "Interactive map of campus itself as a whole"
$afterThis = $d.Paragraphs($afterIdx)
$null = $afterThis.Range.InsertParagraphAfter()

$parkingNew = $d.Paragraphs($afterIdx + 1)
$parkingNew.Range.Text = "Parking"

# ------------------------------------------------------------------
# 3) Plant a fresh "_GoBack" bookmark right before the run of
#    "3-Dimensional campus overview" (an empty bookmark at that
#    point). Word only keeps a single "_GoBack" bookmark, so adding
#    this one automatically removes the old one that used to sit in
#    the "Similar Apps..." paragraph.
# ------------------------------------------------------------------
$threeDIdx = Find-ParaIndex $d "3-Dimensional campus overview"
$threeD = $d.Paragraphs($threeDIdx)
$bmPos = $threeD.Range.Start
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4) Merge the split run "Problems faced by new students: y" +
#    "ou don't know the campus or locations" into a single run.
# ------------------------------------------------------------------
$search4 = "Problems faced by new students: you don" + $apos + "t know the campus or locations"
$null = $d.Content.Find.Execute($search4, $true, $false, $false, $false, $false, $true, 1, $false, $search4, 2)

# ------------------------------------------------------------------
# 5) Merge "Similar Apps to consider for ideas" + ":" into a single
#    run (the old bookmark that used to live between them already
#    moved away in step 3).
# ------------------------------------------------------------------
$search5 = "Similar Apps to consider for ideas:"
$null = $d.Content.Find.Execute($search5, $true, $false, $false, $false, $false, $true, 1, $false, $search5, 2)
